$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: numeric 0, bold, boxed (thin border all sides), centered horizontally / top vertically
$ws.Range("B1").Value = 0
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("B1").Borders.Weight = 2
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4160

# A2: same formatting as B1 - copy B1's style so the style table is reused
# instead of re-deriving it property-by-property (avoids stray cellXfs entries)
$ws.Range("A2").Value = 0
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B2: plain text label
$ws.Range("B2").Value = "disconnected_elements"
